$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above the current row 545 (old "Primera"/"Segunda" pair
# for Terminal La Palmera de La Serena - Acelga). Everything from the old
# row 545 onward shifts down by two rows, preserving formatting (the date
# column D keeps its style, used for the datetime number format).
$ws.Rows("545:546").Insert()

# New row 545: "Primera" quality entry, most recent weekly date (2023-01-02 = 44939)
$ws.Range("A545").Value = 8
$ws.Range("B545").Value = 'Terminal La Palmera de La Serena'
$ws.Range("C545").Value = 'Coquimbo'
$ws.Range("D545").Value = 44939
$ws.Range("E545").Value = 4
$ws.Range("F545").Value = 100112009
$ws.Range("G545").Value = 'Acelga'
$ws.Range("H545").Value = 'Sin especificar'
$ws.Range("I545").Value = 'Primera'
$ws.Range("J545").Value = 2300
$ws.Range("K545").Value = 600
$ws.Range("L545").Value = 700
$ws.Range("M545").Value = 650
$ws.Range("N545").Value = '$/atado 1,5 a 2 kilos'
$ws.Range("O545").Value = 'Provincia del Elquí'
$ws.Range("P545").Value = 325
$ws.Range("Q545").Value = 2
$ws.Range("R545").Value = 'Hortaliza'

# New row 546: "Segunda" quality entry, same date
$ws.Range("A546").Value = 8
$ws.Range("B546").Value = 'Terminal La Palmera de La Serena'
$ws.Range("C546").Value = 'Coquimbo'
$ws.Range("D546").Value = 44939
$ws.Range("E546").Value = 4
$ws.Range("F546").Value = 100112009
$ws.Range("G546").Value = 'Acelga'
$ws.Range("H546").Value = 'Sin especificar'
$ws.Range("I546").Value = 'Segunda'
$ws.Range("J546").Value = 1480
$ws.Range("K546").Value = 500
$ws.Range("L546").Value = 550
$ws.Range("M546").Value = 525
$ws.Range("N546").Value = '$/atado 1,5 a 2 kilos'
$ws.Range("O546").Value = 'Provincia del Elquí'
$ws.Range("P546").Value = 262
$ws.Range("Q546").Value = 2
$ws.Range("R546").Value = 'Hortaliza'
